$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting existing rows 36-38 down to 37-39,
# and copying formatting from the row above (row 35) down - standard Excel
# insert behaviour.
$ws.Rows.Item(36).Insert()

# Fill the newly inserted row 36 with the new weekly data point.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44578
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100114002
$ws.Range("G36").Value = "Camote"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 20000
$ws.Range("N36").Value = "$/malla 20 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 1000
$ws.Range("Q36").Value = 20
$ws.Range("R36").Value = "Hortaliza"
